$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "pontos notáveis - incremento na tabela de ranking"
# Scale the ranking percentage columns (E and F, rows 2-7) up by a factor
# of 100, turning the stored fractional values (0-1) into their
# percentage-point equivalents (0-100), keeping the existing % number format.
$range = $ws.Range("E2:F7")
foreach ($cell in $range.Cells) {
    $cell.Value2 = $cell.Value2 * 100
}
